# Commit: "updated rest assured added header map"
# Adds a new "statusCode" column (with sample value 200) to the
# "createCity" worksheet, and re-applies the column widths that the
# spreadsheet authoring tool recalculated once the extra column existed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createCity")
$ws.Activate()

# New header cell (row 1) and sample value (row 2) for the extra
# "statusCode" field that the header map now exposes.
$ws.Range("H1").Value = "statusCode"
$ws.Range("H2").Value = 200

# Column widths were recomputed across the sheet once the 8th column was
# introduced: A-D and G keep the (new) standard width, E/F shrink very
# slightly, and H gets its own width.
$ws.Columns.Item(1).ColumnWidth = 7.59
$ws.Columns.Item(2).ColumnWidth = 7.59
$ws.Columns.Item(3).ColumnWidth = 7.59
$ws.Columns.Item(4).ColumnWidth = 7.59
$ws.Columns.Item(5).ColumnWidth = 11.25
$ws.Columns.Item(6).ColumnWidth = 18.92
$ws.Columns.Item(7).ColumnWidth = 7.59
$ws.Columns.Item(8).ColumnWidth = 9.59

# Keep the same "active cell" story as before, just shifted from the old
# last column (G2) to the new last column (H2); this also keeps
# createCity as the selected/active sheet and grows the sheet dimension
# to A1:H2.
$ws.Range("H2").Select() | Out-Null
